$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Remove existing hyperlinks on F2:F6 (they will be re-added after data rewrite,
# since row shifting does not carry hyperlink anchors along in this engine).
for ($i = 2; $i -le 6; $i++) {
    $ws.Range("F$i").Hyperlinks.Delete()
}

# Clear old data rows (2-6); full data set will be rewritten below.
$ws.Range("A2:H6").ClearContents()

# Row 2
$ws.Cells.Item(2, 1).Value = "2025-10-19 18:23:32"
$ws.Cells.Item(2, 2).Value = "医療機関向けAIアプリとLINEの連携開発を支援してくださるAIエンジニア募集(AI/バックエンド)"
$ws.Cells.Item(2, 3).Value = "システム開発"
$ws.Cells.Item(2, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(2, 5).Value = "期限情報なし"
$ws.Cells.Item(2, 6).Value = "https://www.lancers.jp/work/detail/5416301"
$ws.Cells.Item(2, 7).Value = 385
$ws.Cells.Item(2, 8).Value = "🔥AI,Ai ◆開発 ◇アプリ"

# Row 3
$ws.Cells.Item(3, 1).Value = "2025-10-19 18:23:32"
$ws.Cells.Item(3, 2).Value = "大企業の業務効率化AIプロジェクトの技術方針策定を支援するAIテックリード募集"
$ws.Cells.Item(3, 3).Value = "システム開発"
$ws.Cells.Item(3, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(3, 5).Value = "期限情報なし"
$ws.Cells.Item(3, 6).Value = "https://www.lancers.jp/work/detail/5416307"
$ws.Cells.Item(3, 7).Value = 378
$ws.Cells.Item(3, 8).Value = "🔥AI,Ai ◆効率化"

# Row 4
$ws.Cells.Item(4, 1).Value = "2025-10-19 18:23:32"
$ws.Cells.Item(4, 2).Value = "Azureでの社内文書検索RAG開発の精度改善を伴走支援してくださるAIエンジニア募集"
$ws.Cells.Item(4, 3).Value = "システム開発"
$ws.Cells.Item(4, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(4, 5).Value = "期限情報なし"
$ws.Cells.Item(4, 6).Value = "https://www.lancers.jp/work/detail/5416305"
$ws.Cells.Item(4, 7).Value = 375
$ws.Cells.Item(4, 8).Value = "🔥AI,Ai ◆開発"

# Row 5
$ws.Cells.Item(5, 1).Value = "2025-10-19 18:23:32"
$ws.Cells.Item(5, 2).Value = "GoogleスプレッドシートとMetaAPIを利用したFXトレード大会ランキングの自動化システム開発"
$ws.Cells.Item(5, 3).Value = "システム開発"
$ws.Cells.Item(5, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(5, 5).Value = "期限情報なし"
$ws.Cells.Item(5, 6).Value = "https://www.lancers.jp/work/detail/5416128"
$ws.Cells.Item(5, 7).Value = 363
$ws.Cells.Item(5, 8).Value = "🔥API ◆開発,システム開発"

# Row 6
$ws.Cells.Item(6, 1).Value = "2025-10-19 18:23:32"
$ws.Cells.Item(6, 2).Value = "Stable Diffusionに詳しいLoRAなどを用いた画像生成AIエンジニア募集"
$ws.Cells.Item(6, 3).Value = "システム開発"
$ws.Cells.Item(6, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(6, 5).Value = "期限情報なし"
$ws.Cells.Item(6, 6).Value = "https://www.lancers.jp/work/detail/5416328"
$ws.Cells.Item(6, 7).Value = 310
$ws.Cells.Item(6, 8).Value = "🔥AI,Ai"

# Row 7
$ws.Cells.Item(7, 1).Value = "2025-10-19 18:23:32"
$ws.Cells.Item(7, 2).Value = "海外仕入れ元サイト→ツールを動かす為のCSVファイルに週1で自動抽出の制作(自動/スクレイピング)"
$ws.Cells.Item(7, 3).Value = "システム開発"
$ws.Cells.Item(7, 4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(7, 5).Value = "期限情報なし"
$ws.Cells.Item(7, 6).Value = "https://www.lancers.jp/work/detail/5251319"
$ws.Cells.Item(7, 7).Value = 135
$ws.Cells.Item(7, 8).Value = "◆ツール,スクレイピング ◇サイト"

# Row 8
$ws.Cells.Item(8, 1).Value = "2025-10-19 18:23:32"
$ws.Cells.Item(8, 2).Value = "イベント出店者管理用ウェブアプリ開発依頼"
$ws.Cells.Item(8, 3).Value = "システム開発"
$ws.Cells.Item(8, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(8, 5).Value = "期限情報なし"
$ws.Cells.Item(8, 6).Value = "https://www.lancers.jp/work/detail/5416005"
$ws.Cells.Item(8, 7).Value = 123
$ws.Cells.Item(8, 8).Value = "◆開発 ◇アプリ"

# Row 9
$ws.Cells.Item(9, 1).Value = "2025-10-19 18:23:32"
$ws.Cells.Item(9, 2).Value = "【恋愛診断】フルスクラッチ開発・運用サポート募集"
$ws.Cells.Item(9, 3).Value = "システム開発"
$ws.Cells.Item(9, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(9, 5).Value = "期限情報なし"
$ws.Cells.Item(9, 6).Value = "https://www.lancers.jp/work/detail/5415986"
$ws.Cells.Item(9, 7).Value = 75
$ws.Cells.Item(9, 8).Value = "◆開発"

# Row 10
$ws.Cells.Item(10, 1).Value = "2025-10-19 18:23:32"
$ws.Cells.Item(10, 2).Value = "【急募】GASを使った顧客管理スプレッドシートの作成・改修依頼"
$ws.Cells.Item(10, 3).Value = "システム開発"
$ws.Cells.Item(10, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(10, 5).Value = "期限情報なし"
$ws.Cells.Item(10, 6).Value = "https://www.lancers.jp/work/detail/5416338"
$ws.Cells.Item(10, 7).Value = 33
$ws.Cells.Item(10, 8).Value = "◇管理"

# Row 11
$ws.Cells.Item(11, 1).Value = "2025-10-19 18:23:32"
$ws.Cells.Item(11, 2).Value = "【高額成功報酬】レガシー基幹システムのバイナリ解析とパッチ作成"
$ws.Cells.Item(11, 3).Value = "システム開発"
$ws.Cells.Item(11, 4).Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Cells.Item(11, 5).Value = "期限情報なし"
$ws.Cells.Item(11, 6).Value = "https://www.lancers.jp/work/detail/5415980"
$ws.Cells.Item(11, 7).Value = 40

# Row 12
$ws.Cells.Item(12, 1).Value = "2025-10-19 18:23:32"
$ws.Cells.Item(12, 2).Value = "【人気調査】どのウォレットや取引所が最も注目されているか?"
$ws.Cells.Item(12, 3).Value = "システム開発"
$ws.Cells.Item(12, 4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(12, 5).Value = "期限情報なし"
$ws.Cells.Item(12, 6).Value = "https://www.lancers.jp/work/detail/5416291"
$ws.Cells.Item(12, 7).Value = 10

# Re-add hyperlinks for F2:F12 pointing at the URL in each cell, matching
# the pre-existing 'Hyperlink' cell style (style index 1) used throughout the sheet.
for ($i = 2; $i -le 12; $i++) {
    $target = $ws.Cells.Item($i, 6).Value2
    $ws.Hyperlinks.Add($ws.Cells.Item($i, 6), $target)
    $ws.Cells.Item($i, 6).Style = "Hyperlink"
}

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
